# Add 2022-Q3 data:
#  1. Insert a new summary row in "总计" for 2022-Q3, shifting existing rows down.
#  2. Insert a new worksheet "2022-Q3" (before "2022-Q2") with the fund holdings detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert new row 2 for 2022-Q3, push others down.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()

# Row 2 (new) lost the A-column border/bold style that the index column
# carries on every other data row - copy it over from row 4 (still has the
# original formatting) before writing values.
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A4").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 1.16

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 1.2

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 3
$wsTotal.Range("D4").Value = 1.3

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q2"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 0.92

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, placed before "2022-Q2".
#    Duplicate the "2022-Q2" sheet so the column styles (bold/bordered
#    header row, bordered index column) come along for free, then overwrite
#    the values and extend it from 4 data rows to 7 data rows.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.ActiveSheet
$wsQ3.Name = "2022-Q3"

# Extend formatting (border on col A, default elsewhere) down to rows 6-8,
# matching row 5's look, before filling in the extra data rows.
$wsQ3.Rows.Item(5).Copy()
$wsQ3.Range("A6:H8").PasteSpecial(-4122)

# Force a cell to stay text even when the literal looks numeric (e.g. fund
# codes like "015466" or ratios like "0.99") - mirrors how these columns
# were authored as inlineStr in the source sheets. Resets to the default
# (un-styled) look, which matches every data cell outside column A here.
function Set-TextCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Header row - plain (non-numeric) text, so it keeps the bordered/bold
# style (s=2) that came along with the worksheet copy.
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Row 2
$wsQ3.Range("A2").Value = 0
Set-TextCell $wsQ3 "B2" "519029"
Set-TextCell $wsQ3 "C2" "华夏稳增混合"
Set-TextCell $wsQ3 "D2" "8.56"
Set-TextCell $wsQ3 "E2" "93.73"
Set-TextCell $wsQ3 "F2" "9.08"
Set-TextCell $wsQ3 "G2" "0.7772"
$wsQ3.Range("H2").Value = 1

# Row 3
$wsQ3.Range("A3").Value = 1
Set-TextCell $wsQ3 "B3" "000480"
Set-TextCell $wsQ3 "C3" "东方红新动力灵活配置混合"
Set-TextCell $wsQ3 "D3" "12.63"
Set-TextCell $wsQ3 "E3" "78.30"
Set-TextCell $wsQ3 "F3" "2.94"
Set-TextCell $wsQ3 "G3" "0.3713"
$wsQ3.Range("H3").Value = 10

# Row 4
$wsQ3.Range("A4").Value = 2
Set-TextCell $wsQ3 "B4" "005997"
Set-TextCell $wsQ3 "C4" "天弘裕利灵活配置混合C"
Set-TextCell $wsQ3 "D4" "0.50"
Set-TextCell $wsQ3 "E4" "44.05"
Set-TextCell $wsQ3 "F4" "1.25"
Set-TextCell $wsQ3 "G4" "0.0062"
$wsQ3.Range("H4").Value = 9

# Row 5
$wsQ3.Range("A5").Value = 3
Set-TextCell $wsQ3 "B5" "015466"
Set-TextCell $wsQ3 "C5" "太平中证1000指数增强A"
Set-TextCell $wsQ3 "D5" "0.37"
Set-TextCell $wsQ3 "E5" "92.23"
Set-TextCell $wsQ3 "F5" "0.99"
Set-TextCell $wsQ3 "G5" "0.0037"
$wsQ3.Range("H5").Value = 8

# Row 6
$wsQ3.Range("A6").Value = 4
Set-TextCell $wsQ3 "B6" "002020"
Set-TextCell $wsQ3 "C6" "国都创新驱动灵活配置混合"
Set-TextCell $wsQ3 "D6" "0.12"
Set-TextCell $wsQ3 "E6" "65.45"
Set-TextCell $wsQ3 "F6" "2.62"
Set-TextCell $wsQ3 "G6" "0.0031"
$wsQ3.Range("H6").Value = 10

# Row 7
$wsQ3.Range("A7").Value = 5
Set-TextCell $wsQ3 "B7" "002388"
Set-TextCell $wsQ3 "C7" "天弘裕利灵活配置混合A"
Set-TextCell $wsQ3 "D7" "0.10"
Set-TextCell $wsQ3 "E7" "44.05"
Set-TextCell $wsQ3 "F7" "1.25"
Set-TextCell $wsQ3 "G7" "0.0012"
$wsQ3.Range("H7").Value = 9

# Row 8
$wsQ3.Range("A8").Value = 6
Set-TextCell $wsQ3 "B8" "015467"
Set-TextCell $wsQ3 "C8" "太平中证1000指数增强C"
Set-TextCell $wsQ3 "D8" "0.02"
Set-TextCell $wsQ3 "E8" "92.23"
Set-TextCell $wsQ3 "F8" "0.99"
Set-TextCell $wsQ3 "G8" "0.0002"
$wsQ3.Range("H8").Value = 8

Write-Output "done"
